$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Date and Time Clocked-In for the new entry on row 7.
# A leading apostrophe forces the values to be stored as text rather
# than being auto-converted into date/time serial numbers.
$ws.Range("A7").Value = "'2026-01-24"
$ws.Range("B7").Value = "22:59:38"

# Copy the formatting from the neighboring cell on the same row so the
# new cells match the existing row styling (font/border).
$ws.Range("C7").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
